# 13/03 - EOD Commit
#
# Append three new sets of profile credentials (Approver, Customer
# Support, Sales Support - each a user/pass pair) to the bottom of the
# "Environment_DirectSales" sheet, and leave that sheet active/selected
# (previously "GeneralVariables" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environment_DirectSales")

$newRows = @(
    @("approverProfileUser",        "andre.esteves.ext@proximus.com.aprv"),
    @("approverProfilePass",        "Inno6677!"),
    @("costumerSupportProfileUser", "andre.esteves.ext@proximus.com.csupp"),
    @("costumerSupportProfilePass", "Inno6677!"),
    @("salesSupportProfileUser",    "andre.esteves.ext@proximus.com.ssupp"),
    @("salesSupportProfilePass",    "Win6677!")
)

$row = 15
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Make "Environment_DirectSales" the active/selected sheet (this also
# clears tabSelected on the formerly-active "GeneralVariables" sheet,
# which keeps its own previously-saved selection).
$ws.Activate() | Out-Null
$ws.Range("B27").Select() | Out-Null

# Best-effort: nudge the saved window position, matching the small
# vertical shift recorded for the workbook view (yWindow -120 -> -90).
try {
    $excel.ActiveWindow.Top = -90
} catch {
}
